$wb = $excel.ActiveWorkbook
$wsBinary = $wb.Worksheets.Item("Binary")
$wsWeights = $wb.Worksheets.Item("Weights")

# --- Header row additions: columns L:O (12:15) = 11,12,13,14 as text, matching existing header style ---
$newHeaders = @("11","12","13","14")
foreach ($ws in @($wsBinary, $wsWeights)) {
    $hdrRange = $ws.Range("L1:O1")
    $hdrRange.NumberFormat = "@"
    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $ws.Cells.Item(1, 12 + $i).Value = $newHeaders[$i]
    }
    $srcHeader = $ws.Range("K1")
    $srcHeader.Copy()
    $hdrRange.PasteSpecial(-4122)  # xlPasteFormats: reuse bold/border/center style
}

# --- Binary sheet: full 15x15 adjacency data (rows 2-16, cols A-O) ---
$binaryData = @{
    2 = @(0,1,1,1,1,0,0,0,0,0,0,0,0,0,0)
    3 = @(1,0,1,1,1,0,0,0,0,0,0,0,0,0,0)
    4 = @(1,1,0,1,1,0,0,0,0,0,0,0,0,0,0)
    5 = @(1,1,1,0,1,0,0,0,0,0,0,0,0,0,0)
    6 = @(1,1,1,1,0,0,0,0,0,0,0,0,0,0,0)
    7 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    8 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    9 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    10 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    11 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    12 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    13 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    14 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    15 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    16 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
}
foreach ($r in $binaryData.Keys) {
    $rowVals = $binaryData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $wsBinary.Cells.Item([int]$r, $c + 1).Value = $rowVals[$c]
    }
}

# --- Weights sheet: full 15x15 weight data (rows 2-16, cols A-O) ---
$weightsData = @{
    2 = @(0,66,95,45,81,0,0,0,0,0,0,0,0,0,0)
    3 = @(66,0,83,51,36,0,0,0,0,0,0,0,0,0,0)
    4 = @(95,83,0,68,53,0,0,0,0,0,0,0,0,0,0)
    5 = @(45,51,68,0,11,0,0,0,0,0,0,0,0,0,0)
    6 = @(81,36,53,11,0,0,0,0,0,0,0,0,0,0,0)
    7 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    8 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    9 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    10 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    11 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    12 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    13 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    14 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    15 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
    16 = @(0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)
}
foreach ($r in $weightsData.Keys) {
    $rowVals = $weightsData[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $wsWeights.Cells.Item([int]$r, $c + 1).Value = $rowVals[$c]
    }
}
